$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.050905999999999
$ws.Range("H2").Value = 24.152718
$ws.Range("I2").Value = 0.1888708516018927
$ws.Range("J2").Value = 0.1888708516018927
$ws.Range("M2").Value = 179.7005413333333
$ws.Range("N2").Value = 539.101624
$ws.Range("O2").Value = 0.7012656334041908
$ws.Range("P2").Value = 0.7012656334041907
$ws.Range("Q2").Value = 1446.752166423781
$ws.Range("R2").Value = 13020.76949781403
$ws.Range("S2").Value = 0.1324486373801902
$ws.Range("T2").Value = 0.1324486373801902
$ws.Range("G3").Value = 8.050905999999999
$ws.Range("H3").Value = 24.152718
$ws.Range("I3").Value = 0.1888708516018927
$ws.Range("J3").Value = 0.1888708516018927
$ws.Range("O3").Value = 0.05908927597267952
$ws.Range("P3").Value = 0.05908927597267952
$ws.Range("Q3").Value = 121.90464490738
$ws.Range("R3").Value = 1097.14180416642
$ws.Range("S3").Value = 0.01116024187349924
$ws.Range("T3").Value = 0.01116024187349924
$ws.Range("G4").Value = 8.050905999999999
$ws.Range("H4").Value = 24.152718
$ws.Range("I4").Value = 0.1888708516018927
$ws.Range("J4").Value = 0.1888708516018927
$ws.Range("M4").Value = 36.14947766666667
$ws.Range("N4").Value = 108.448433
$ws.Range("O4").Value = 0.1410701724382803
$ws.Range("P4").Value = 0.1410701724382803
$ws.Range("Q4").Value = 291.0360466434327
$ws.Range("R4").Value = 2619.324419790894
$ws.Range("S4").Value = 0.02664404360404385
$ws.Range("T4").Value = 0.02664404360404385
$ws.Range("G5").Value = 8.050905999999999
$ws.Range("H5").Value = 24.152718
$ws.Range("I5").Value = 0.1888708516018927
$ws.Range("J5").Value = 0.1888708516018927
$ws.Range("M5").Value = 25.25999466666667
$ws.Range("N5").Value = 75.779984
$ws.Range("O5").Value = 0.09857491818484938
$ws.Range("P5").Value = 0.09857491818484938
$ws.Range("Q5").Value = 203.3658426218347
$ws.Range("R5").Value = 1830.292583596512
$ws.Range("S5").Value = 0.0186179287441594
$ws.Range("T5").Value = 0.0186179287441594
$ws.Range("I6").Value = 0.2904749299149038
$ws.Range("J6").Value = 0.2904749299149038
$ws.Range("M6").Value = 179.7005413333333
$ws.Range("N6").Value = 539.101624
$ws.Range("O6").Value = 0.7012656334041908
$ws.Range("P6").Value = 0.7012656334041907
$ws.Range("Q6").Value = 2225.040182653425
$ws.Range("R6").Value = 20025.36164388082
$ws.Range("S6").Value = 0.203700085714813
$ws.Range("T6").Value = 0.2037000857148129
$ws.Range("I7").Value = 0.2904749299149038
$ws.Range("J7").Value = 0.2904749299149038
$ws.Range("O7").Value = 0.05908927597267952
$ws.Range("P7").Value = 0.05908927597267952
$ws.Range("S7").Value = 0.01716395329688649
$ws.Range("T7").Value = 0.01716395329688649
$ws.Range("I8").Value = 0.2904749299149038
$ws.Range("J8").Value = 0.2904749299149038
$ws.Range("M8").Value = 36.14947766666667
$ws.Range("N8").Value = 108.448433
$ws.Range("O8").Value = 0.1410701724382803
$ws.Range("P8").Value = 0.1410701724382803
$ws.Range("Q8").Value = 447.6004345533148
$ws.Range("R8").Value = 4028.403910979833
$ws.Range("S8").Value = 0.04097734845209287
$ws.Range("T8").Value = 0.04097734845209287
$ws.Range("I9").Value = 0.2904749299149038
$ws.Range("J9").Value = 0.2904749299149038
$ws.Range("M9").Value = 25.25999466666667
$ws.Range("N9").Value = 75.779984
$ws.Range("O9").Value = 0.09857491818484938
$ws.Range("P9").Value = 0.09857491818484938
$ws.Range("Q9").Value = 312.7675783830205
$ws.Range("R9").Value = 2814.908205447184
$ws.Range("S9").Value = 0.0286335424511115
$ws.Range("T9").Value = 0.0286335424511115
$ws.Range("G10").Value = 8.230170000000001
$ws.Range("H10").Value = 24.69051
$ws.Range("I10").Value = 0.1930763092661061
$ws.Range("J10").Value = 0.1930763092661061
$ws.Range("M10").Value = 179.7005413333333
$ws.Range("N10").Value = 539.101624
$ws.Range("O10").Value = 0.7012656334041908
$ws.Range("P10").Value = 0.7012656334041907
$ws.Range("Q10").Value = 1478.96600426536
$ws.Range("R10").Value = 13310.69403838824
$ws.Range("S10").Value = 0.1353977803128393
$ws.Range("T10").Value = 0.1353977803128393
$ws.Range("G11").Value = 8.230170000000001
$ws.Range("H11").Value = 24.69051
$ws.Range("I11").Value = 0.1930763092661061
$ws.Range("J11").Value = 0.1930763092661061
$ws.Range("O11").Value = 0.05908927597267952
$ws.Range("P11").Value = 0.05908927597267952
$ws.Range("Q11").Value = 124.6190119941
$ws.Range("R11").Value = 1121.5711079469
$ws.Range("S11").Value = 0.01140873932201136
$ws.Range("T11").Value = 0.01140873932201136
$ws.Range("G12").Value = 8.230170000000001
$ws.Range("H12").Value = 24.69051
$ws.Range("I12").Value = 0.1930763092661061
$ws.Range("J12").Value = 0.1930763092661061
$ws.Range("M12").Value = 36.14947766666667
$ws.Range("N12").Value = 108.448433
$ws.Range("O12").Value = 0.1410701724382803
$ws.Range("P12").Value = 0.1410701724382803
$ws.Range("Q12").Value = 297.51634660787
$ws.Range("R12").Value = 2677.647119470831
$ws.Range("S12").Value = 0.02723730824191632
$ws.Range("T12").Value = 0.02723730824191632
$ws.Range("G13").Value = 8.230170000000001
$ws.Range("H13").Value = 24.69051
$ws.Range("I13").Value = 0.1930763092661061
$ws.Range("J13").Value = 0.1930763092661061
$ws.Range("M13").Value = 25.25999466666667
$ws.Range("N13").Value = 75.779984
$ws.Range("O13").Value = 0.09857491818484938
$ws.Range("P13").Value = 0.09857491818484938
$ws.Range("Q13").Value = 207.89405030576
$ws.Range("R13").Value = 1871.04645275184
$ws.Range("S13").Value = 0.01903248138933909
$ws.Range("T13").Value = 0.01903248138933909
$ws.Range("G14").Value = 13.96350433333333
$ws.Range("H14").Value = 41.890513
$ws.Range("I14").Value = 0.3275779092170975
$ws.Range("J14").Value = 0.3275779092170975
$ws.Range("M14").Value = 179.7005413333333
$ws.Range("N14").Value = 539.101624
$ws.Range("O14").Value = 0.7012656334041908
$ws.Range("P14").Value = 0.7012656334041907
$ws.Range("Q14").Value = 2509.249287610346
$ws.Range("R14").Value = 22583.24358849311
$ws.Range("S14").Value = 0.2297191299963484
$ws.Range("T14").Value = 0.2297191299963483
$ws.Range("G15").Value = 13.96350433333333
$ws.Range("H15").Value = 41.890513
$ws.Range("I15").Value = 0.3275779092170975
$ws.Range("J15").Value = 0.3275779092170975
$ws.Range("O15").Value = 0.05908927597267952
$ws.Range("P15").Value = 0.05908927597267952
$ws.Range("Q15").Value = 211.4316124691633
$ws.Range("R15").Value = 1902.88451222247
$ws.Range("S15").Value = 0.01935634148028243
$ws.Range("T15").Value = 0.01935634148028243
$ws.Range("G16").Value = 13.96350433333333
$ws.Range("H16").Value = 41.890513
$ws.Range("I16").Value = 0.3275779092170975
$ws.Range("J16").Value = 0.3275779092170975
$ws.Range("M16").Value = 36.14947766666667
$ws.Range("N16").Value = 108.448433
$ws.Range("O16").Value = 0.1410701724382803
$ws.Range("P16").Value = 0.1410701724382803
$ws.Range("Q16").Value = 504.7733880462366
$ws.Range("R16").Value = 4542.960492416129
$ws.Range("S16").Value = 0.04621147214022726
$ws.Range("T16").Value = 0.04621147214022726
$ws.Range("G17").Value = 13.96350433333333
$ws.Range("H17").Value = 41.890513
$ws.Range("I17").Value = 0.3275779092170975
$ws.Range("J17").Value = 0.3275779092170975
$ws.Range("M17").Value = 25.25999466666667
$ws.Range("N17").Value = 75.779984
$ws.Range("O17").Value = 0.09857491818484938
$ws.Range("P17").Value = 0.09857491818484938
$ws.Range("Q17").Value = 352.7180449879769
$ws.Range("R17").Value = 3174.462404891792
$ws.Range("S17").Value = 0.0322909656002394
$ws.Range("T17").Value = 0.0322909656002394
